$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 41; existing rows 41..150 shift down to 42..151.
$ws.Rows.Item(41).Insert()

# Copy the number format (date style) that row 42 (the row that used to be row 41,
# now shifted down) has for column D, so the newly inserted row's D cell matches.
$ws.Range("D41").NumberFormat = $ws.Range("D42").NumberFormat

# Populate the new row 41 with its data.
$ws.Range("A41").Value = 7
$ws.Range("B41").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C41").Value = "Ñuble"
$ws.Range("D41").Value = 44526
$ws.Range("E41").Value = 16
$ws.Range("F41").Value = 100112032
$ws.Range("G41").Value = "Zapallo italiano"
$ws.Range("H41").Value = "Sin especificar"
$ws.Range("I41").Value = "Primera"
$ws.Range("J41").Value = 60
$ws.Range("K41").Value = 8000
$ws.Range("L41").Value = 9000
$ws.Range("M41").Value = 8500
$ws.Range("N41").Value = "$/caja 60 unidades"
$ws.Range("O41").Value = "Región del Maule"
$ws.Range("P41").Value = 142
$ws.Range("Q41").Value = 60
$ws.Range("R41").Value = "Hortaliza"
